# Replace CSV mapping files:
#  - Remove repeated symbol (row 21 "Obliqua1 middle", duplicate of row 5's neume.obliquamiddle1 entry)
#  - Change width for obliques (rows 5-7) to 2

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Change width (column G) for the three "oblique" rows from 1 to 2
$ws.Range("G5").Value = 2
$ws.Range("G6").Value = 2
$ws.Range("G7").Value = 2

# Remove the repeated "Obliqua1 middle" row (row 21), shifting rows 22-25 up
$ws.Rows("21:21").Delete()

# Recompute the wrapped-text row heights for the rows that moved / were affected
$ws.Rows("19:19").RowHeight = 68
$ws.Rows("20:20").RowHeight = 68
$ws.Rows("21:21").RowHeight = 68
$ws.Rows("22:22").RowHeight = 85
$ws.Rows("23:23").RowHeight = 85
$ws.Rows("24:24").RowHeight = 85

# Reflect final selection position left by the edit
$ws.Range("H22").Select() | Out-Null
